# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders/updates the worker-period rows (16-23) on Hoja1 so that the two
# workers (YESSICA DEL CARMEN GARCIA CARREAZO / LUC ENGETSCHWILER) are
# interleaved in ascending period order (2405..2408), and refreshes the
# "Valor Mora" (F) / "Salario Basico" (G) amounts per the new account
# statement data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Tipo Doc Trabajador (B), N Doc Trabajador (C), Nombre Trabajador (D),
# Periodo Mora (E), Valor Mora (F), Salario Basico (G)
$rows = @(
    @("CC", "1047418337", "YESSICA DEL CARMEN GARCIA CARREAZO", "2405", 72000, 1800000),
    @("CC", "2000012246", "LUC ENGETSCHWILER",                  "2405", 52000, 1300000),
    @("CC", "1047418337", "YESSICA DEL CARMEN GARCIA CARREAZO", "2406", 72000, 1800000),
    @("CC", "2000012246", "LUC ENGETSCHWILER",                  "2406", 52000, 1300000),
    @("CC", "1047418337", "YESSICA DEL CARMEN GARCIA CARREAZO", "2407", 72000, 1800000),
    @("CC", "2000012246", "LUC ENGETSCHWILER",                  "2407", 52000, 1300000),
    @("CC", "1047418337", "YESSICA DEL CARMEN GARCIA CARREAZO", "2408", 50400, 1800000),
    @("CC", "2000012246", "LUC ENGETSCHWILER",                  "2408", 36400, 1300000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
}
